$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.472.99"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "3.006.67"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'583.84"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").Value = "'146.29"
$ws.Range("E6").Value = "  -6.37%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").Value = "3.004.41"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  -5.45%  "
$ws.Range("D11").Value = "'5.71"
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").Value = "'0.443"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "'0.0000228"
$ws.Range("E13").Value = "  -4.45%  "
$ws.Range("D14").Value = "'34.63"
$ws.Range("E14").Value = "  -6.13%  "
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "3.499.56"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "'7.04"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "62.427.15"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "3.002.50"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("D20").Value = "'459.60"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "'13.89"
$ws.Range("E21").Value = "  -3.85%  "
$ws.Range("D22").Value = "'0.681"
$ws.Range("E22").Value = "  -4.34%  "
$ws.Range("D23").Value = "'7.32"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("D24").Value = "'80.03"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "  -7.02%  "
$ws.Range("D26").Value = "'12.27"
$ws.Range("E26").Value = "  -4.70%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'9.99"
$ws.Range("E28").Value = "  -5.69%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'7.21"
$ws.Range("E30").Value = "  -3.99%  "
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("D33").Value = "'27.02"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("D35").Value = "'1.02"
$ws.Range("E35").Value = "  -3.08%  "
$ws.Range("D36").Value = "0.0₃0790"
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("D38").Value = "'2.12"
$ws.Range("E38").Value = "  -5.97%  "
$ws.Range("D39").Value = "'50.21"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").Value = "'8.96"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").Value = "'2.94"
$ws.Range("E41").Value = "  -10.17%  "
$ws.Range("D42").Value = "'411.21"
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("D43").Value = "'0.280"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "2.774.64"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").Value = "'39.19"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "'0.0352"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "'127.26"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "'23.75"
$ws.Range("E51").Value = "  -7.95%  "